$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price value is a plain decimal number (e.g. "1.002") must be
# forced to Text format first, otherwise Excel auto-converts the literal string
# into a numeric value and loses formatting like trailing/leading zeros.
$textForcedCells = @(
    "D4",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D14",
    "D15",
    "D16",
    "D18",
    "D20",
    "D22",
    "D23",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D34",
    "D35",
    "D37",
    "D38",
    "D39",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.922.12"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.816.71"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "0.4671"
$ws.Range("E7").Value = "  +1.20%  "
$ws.Range("D8").Value = "0.3688"
$ws.Range("E8").Value = "  -1.55%  "
$ws.Range("D9").Value = "0.07363"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("D10").Value = "0.8712"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").Value = "20.39"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").Value = "1.828.01"
$ws.Range("E12").Value = "  +5.28%  "
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").Value = "6.518"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "0.07070"
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").Value = "91.60"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").Value = "0.000008703"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").Value = "14.73"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").Value = "26.933.30"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "5.323"
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").Value = "10.62"
$ws.Range("E23").Value = "  -0.82%  "
$ws.Range("D24").Value = "2.025.04"
$ws.Range("E24").Value = "  +3.48%  "
$ws.Range("D25").Value = "1.895"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("D26").Value = "150.48"
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").Value = "2.171"
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("D28").Value = "18.35"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "5.337"
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("D30").Value = "116.24"
$ws.Range("E30").Value = "  +1.30%  "
$ws.Range("D31").Value = "0.08956"
$ws.Range("E31").Value = "  +0.76%  "
$ws.Range("D32").Value = "0.7681"
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("D34").Value = "4.507"
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("D35").Value = "2.924"
$ws.Range("E35").Value = "  +1.30%  "
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("D37").Value = "1.085"
$ws.Range("E37").Value = "  -3.15%  "
$ws.Range("D38").Value = "0.01965"
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("D39").Value = "0.05292"
$ws.Range("E39").Value = "  +1.11%  "
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("D41").Value = "7.262"
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("D42").Value = "0.5321"
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("D43").Value = "2.348"
$ws.Range("E43").Value = "  -5.79%  "
$ws.Range("D44").Value = "0.1657"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Value = "8.446"
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("D46").Value = "0.4926"
$ws.Range("E46").Value = "  -2.49%  "
$ws.Range("D47").Value = "10.46"
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").Value = "1.672"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").Value = "103.82"
$ws.Range("E50").Value = "  -0.51%  "
$ws.Range("D51").Value = "0.06301"
$ws.Range("E51").Value = "  -0.23%  "
